$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 439
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value()
    $cell.Value = $cur.AddDays(1)
}
